$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "fog=true/fog_visible=85" note in B6 is no longer needed -> clear it
# (this also drops the now-unused shared string from sharedStrings.xml)
$ws.Range("B6").ClearContents()

# Column B needs to be noticeably wider to fit the (now longer) notes column
$ws.Columns.Item(2).ColumnWidth = 51.83

# Reflect where the user was working when they saved: zoomed in on B6
$ws.Range("B6").Select()
$excel.ActiveWindow.Zoom = 115
